$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# --- Row 35/36: coins re-ranked (Bittensor now ahead of Hedera) ---
Set-TextValue "B35" "Bittensor"
Set-TextValue "C35" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D35" "619.42"
Set-TextValue "E35" "  +1.10%  "
Set-TextValue "B36" "Hedera"
Set-TextValue "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.118"
Set-TextValue "E36" "  +1.74%  "

# --- Row 39/40: coins re-ranked (FirstDigitalUSD now ahead of PEPE) ---
Set-TextValue "B39" "FirstDigitalUSD"
Set-TextValue "C39" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.00%  "
Set-TextValue "B40" "PEPE"
Set-TextValue "C40" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D40" "0.0₃0776"
Set-TextValue "E40" "  -12.58%  "

# --- Price / Volume(1h) refresh for remaining rows ---
Set-TextValue "D2" "67.983.55"
Set-TextValue "E2" "  +0.03%  "
Set-TextValue "D3" "3.670.81"
Set-TextValue "E3" "  -0.98%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "601.02"
Set-TextValue "E5" "  +3.48%  "
Set-TextValue "D6" "192.36"
Set-TextValue "E6" "  +9.59%  "
Set-TextValue "D7" "0.624"
Set-TextValue "E7" "  -0.26%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  +0.16%  "
Set-TextValue "D9" "0.707"
Set-TextValue "E9" "  +0.30%  "
Set-TextValue "D10" "58.15"
Set-TextValue "E10" "  +12.45%  "
Set-TextValue "D11" "0.153"
Set-TextValue "E11" "  -5.49%  "
Set-TextValue "E12" "  -4.92%  "
Set-TextValue "D13" "10.21"
Set-TextValue "E13" "  -1.89%  "
Set-TextValue "D14" "4.260.70"
Set-TextValue "E14" "  -0.75%  "
Set-TextValue "D15" "3.675.13"
Set-TextValue "E15" "  -0.58%  "
Set-TextValue "E16" "  +0.84%  "
Set-TextValue "D17" "19.00"
Set-TextValue "E17" "  -1.84%  "
Set-TextValue "E18" "  +0.66%  "
Set-TextValue "D19" "67.845.41"
Set-TextValue "E19" "  +0.08%  "
Set-TextValue "D20" "12.53"
Set-TextValue "E20" "  -2.42%  "
Set-TextValue "D21" "400.77"
Set-TextValue "E21" "  -1.11%  "
Set-TextValue "E22" "  -0.64%  "
Set-TextValue "D23" "88.19"
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "11.49"
Set-TextValue "E24" "  +6.45%  "
Set-TextValue "E25" "  -2.10%  "
Set-TextValue "E26" "  -1.04%  "
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "E28" "  -2.71%  "
Set-TextValue "D29" "9.31"
Set-TextValue "E29" "  -1.58%  "
Set-TextValue "D30" "31.98"
Set-TextValue "E30" "  -1.64%  "
Set-TextValue "E31" "  +2.40%  "
Set-TextValue "D32" "45.78"
Set-TextValue "E32" "  +6.38%  "
Set-TextValue "E33" "  -0.33%  "
Set-TextValue "D34" "67.11"
Set-TextValue "E34" "  +3.49%  "
Set-TextValue "E37" "  -0.08%  "
Set-TextValue "D38" "0.399"
Set-TextValue "E38" "  +0.79%  "
Set-TextValue "E41" "  -0.43%  "
Set-TextValue "E42" "  -2.24%  "
Set-TextValue "E43" "  -1.34%  "
Set-TextValue "E44" "  -7.03%  "
Set-TextValue "D45" "2.838.66"
Set-TextValue "E45" "  +0.95%  "
Set-TextValue "D46" "0.136"
Set-TextValue "E46" "  +2.25%  "
Set-TextValue "D47" "3.22"
Set-TextValue "E47" "  +5.89%  "
Set-TextValue "D48" "8.97"
Set-TextValue "E48" "  -2.64%  "
Set-TextValue "D49" "144.61"
Set-TextValue "E49" "  +4.33%  "
Set-TextValue "D50" "2.64"
Set-TextValue "E50" "  -1.65%  "
Set-TextValue "E51" "  -12.70%  "
